# fix: fixed formatting when scrapping floating point numbers
#
# 1) A handful of "Razon social"/"Nombre Fantasia" values used a comma (",")
#    as a separator between multiple people/companies, which is ambiguous
#    with the comma sometimes also appearing inside abbreviations like
#    "S.H." (Sociedad de Hecho). These are corrected to use a period (".")
#    as the separator, and the "S.H." abbreviation is normalized to "SH".
#
# 2) The "Importe" column (H) was scraped with Spanish/Argentine number
#    formatting (period as thousands separator, comma as decimal
#    separator), e.g. "58.960,00". These values are stored as plain text
#    and need to be reformatted to a plain decimal representation (no
#    thousands separator, period as decimal separator), e.g. "58960.00".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix proveedor name fields that used commas as list separators ---
$ws.Range("E111").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E193").Value = "EDICIONES NATIVA SH DE ESCOBAR JORGE. MARTINEZ ALFREDO. PIZIGHINI CARLOS L Y R"
$ws.Range("E222").Value = "OLIVERA. FLORENCIO"
$ws.Range("F222").Value = "OLIVERA. FLORENCIO"
$ws.Range("E226").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E240").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("F310").Value = "CLERICE. MIGUEL ANGEL"

# --- 2) Reformat the "Importe" column values from "1.234,56" to "1234.56" ---
# Force the column to Text format first so Excel keeps storing these as
# plain text strings (matching the original file) instead of coercing
# the now-plain-decimal-looking values into real numbers.
$ws.Range("H2:H330").NumberFormat = "@"
for ($r = 2; $r -le 330; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Value()
    $new = $old.Replace(".", "").Replace(",", ".")
    $cell.Value = $new
}
